$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F10").Value = -3
$ws.Range("F11").Value = -7
$ws.Range("F14").Value = -13
$ws.Range("F18").Value = -2
$ws.Range("F21").Value = 1
$ws.Range("F22").Value = -3
$ws.Range("F24").Value = -6
$ws.Range("F26").Value = -5
$ws.Range("F28").Value = 7
$ws.Range("F29").Value = -1
$ws.Range("F30").Value = -1
$ws.Range("F32").Value = -8
$ws.Range("F34").Value = 1
$ws.Range("F35").Value = -2
$ws.Range("F38").Value = -2
$ws.Range("F39").Value = 12
$ws.Range("F40").Value = -2
$ws.Range("F41").Value = -4
$ws.Range("F44").Value = 1
$ws.Range("F45").Value = 4
